$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# This shared string is used on the Overview sheet (columns E/F = zh-cn/de-de
# status) as well as on the per-locale sheets (zh-cn, de-de) in the "Status"
# column. Using Find/Replace across every worksheet keeps it in sync no
# matter the sheet layout.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- Narrow the status columns ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
